$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndexes = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

$data = @(
    @(1.02, 1.026080828941316, 1.029432941799755, 1.029709536245137, 1.035662765503488, 1.030734286204386, 1.031246262710298, 1.032246867966284, 1.032522659778654, 1.038458726611105, 1.014325558436712),
    @(1.02, 1.027406892461928, 1.03038500118079, 1.030983780672548, 1.037124521082735, 1.031000977895484, 1.032209994537634, 1.033006684220129, 1.033603853017924, 1.039728191114787, 1.014655009853402),
    @(1.02, 1.028264018477896, 1.031000082699094, 1.031807763777968, 1.038069597556168, 1.031171759026243, 1.032832248903379, 1.033496757783729, 1.034302376647009, 1.040548339696167, 1.014867467759454),
    @(1.02, 1.028624137326776, 1.031258434404268, 1.032154041265391, 1.038466727102832, 1.031243128766361, 1.033093525914171, 1.033702409161418, 1.034595780699788, 1.040892828547104, 1.014956613837086),
    @(1.02, 1.028684590281532, 1.031301799457286, 1.032212175560042, 1.038533396468339, 1.031255087069286, 1.033137376875558, 1.033736916977276, 1.034645029661135, 1.040950652192914, 1.014971571842474),
    @(1.02, 1.028268831250128, 1.031003535702424, 1.031812391241135, 1.038074904723501, 1.031172714347431, 1.032835741348035, 1.033499507180297, 1.034306298125142, 1.040552943952261, 1.014868659605334),
    @(1.02, 1.026529172833983, 1.029754895865943, 1.030140286883499, 1.036156936838686, 1.030824786351502, 1.031572240707298, 1.032503979633974, 1.032888279936165, 1.038888015968032, 1.014437047572193),
    @(1.02, 1.023456371203236, 1.027547139235124, 1.027189530032637, 1.03277103602332, 1.030197968133543, 1.029335350755811, 1.030737540158177, 1.030381119127288, 1.035944183200685, 1.013670940994407),
    @(1.02, 1.021402630717923, 1.026070111814781, 1.025219226595167, 1.03050925897561, 1.029770798453171, 1.027836869735381, 1.029551567955509, 1.028703791045571, 1.033974590714848, 1.013156410052118),
    @(1.02, 1.02051203687685, 1.025429277217638, 1.024365265469818, 1.029528737784798, 1.029583611119971, 1.02718625266947, 1.029036015710176, 1.027976043690221, 1.03311999834682, 1.0129326988663),
    @(1.02, 1.020181027924973, 1.025191048543188, 1.024047940248948, 1.029164348036651, 1.029513746562614, 1.026944315357003, 1.028844210409825, 1.027705503248165, 1.032802296460576, 1.012849463727906),
    @(1.02, 1.020252039737502, 1.025242158238827, 1.024116013367368, 1.02924251915189, 1.029528747921034, 1.026996223976782, 1.028885367204195, 1.027763545237705, 1.032870456827376, 1.012867324250622),
    @(1.02, 1.020484679743449, 1.025409589140001, 1.024339037867246, 1.029498620936139, 1.02957784293227, 1.027166259569935, 1.029020167288834, 1.027953685300524, 1.033093742516097, 1.012925821463369),
    @(1.02, 1.020627989888317, 1.025512723022136, 1.024476433776188, 1.029656389575331, 1.029608047579119, 1.0272709883233, 1.029103181400095, 1.028070807232188, 1.033231280476875, 1.012961845091423),
    @(1.02, 1.02146170838683, 1.026112614835761, 1.025275883773789, 1.030574307992302, 1.029783174574582, 1.027880011521332, 1.029585740678763, 1.028752058253132, 1.034031269907131, 1.013171237624277),
    @(1.02, 1.021984322965797, 1.026488568416201, 1.025777138247396, 1.031149779520702, 1.029892431821171, 1.028261560347926, 1.029887894862435, 1.029178996894456, 1.034532610552387, 1.013302337944791),
    @(1.02, 1.022289029144944, 1.026707733261433, 1.026069433992122, 1.031485331169291, 1.029955945612194, 1.028483941309394, 1.030063941738053, 1.029427883060864, 1.034824866025907, 1.01337871825806),
    @(1.02, 1.022392904870889, 1.026782442120987, 1.026169086285939, 1.031599726958179, 1.029977565895747, 1.0285597387537, 1.030123936264366, 1.029512723209502, 1.03492448926563, 1.013404747033507),
    @(1.02, 1.021928264450479, 1.026448244850862, 1.025723366439714, 1.031088048419837, 1.029880731709211, 1.028220641385284, 1.029855496723666, 1.029133204923984, 1.03447883887178, 1.013288281262304),
    @(1.02, 1.020416178703898, 1.025360290297604, 1.024273366228967, 1.029423210381289, 1.029563394929127, 1.027116195788401, 1.028980480501613, 1.027897699972962, 1.033027997898366, 1.012908599334989),
    @(1.02, 1.019464294442566, 1.024675125130621, 1.023360961922453, 1.028375412836066, 1.029361934720495, 1.026420227692544, 1.02842854975687, 1.027119598495176, 1.032114241730421, 1.012669074311691),
    @(1.02, 1.019969019721255, 1.025038451814139, 1.023844715679378, 1.028930971713094, 1.029468916727118, 1.026789322698519, 1.028721307781661, 1.027532208589167, 1.032598790538993, 1.012796127678123),
    @(1.02, 1.021953595269682, 1.026466465712765, 1.025747663844249, 1.031115942390134, 1.02988601914769, 1.028239131428248, 1.02987013664898, 1.02915389677643, 1.034503136498494, 1.013294633144789),
    @(1.02, 1.024251659800803, 1.028118800783787, 1.027952905479923, 1.033647142116169, 1.030361648579983, 1.029914898330228, 1.031195667189274, 1.031030301898483, 1.036706452626415, 1.013869661983548)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $colIndexes.Length; $c++) {
        $ws.Cells.Item($excelRow, $colIndexes[$c]).Value = $rowValues[$c]
    }
}
